# Fruta / hortaliza, semanal
# Insert a new week's worth of observations (3 quality grades for
# "Frutilla" dated 2022-10-21) right before the existing block of rows
# that start at row 29, pushing the remaining data rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 29..31; everything from the old row 29 onward
# (old rows 29-39) shifts down to 32-42, keeping per-column formatting
# (e.g. the date style on column D) intact.
$ws.Rows("29:31").Insert()

# New rows: Row, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg
$newRows = @(
    @(29, 'Especial', 160, 6500, 7000, 6750, 2250),
    @(30, 'Primera',  160, 5500, 6000, 5750, 1917),
    @(31, 'Segunda',  200, 4500, 5000, 4750, 1583)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = 'Agrícola del Norte S.A. de Arica'
    $ws.Cells.Item($row, 3).Value = 'Arica y Parinacota'
    $ws.Cells.Item($row, 4).Value = 44855
    $ws.Cells.Item($row, 5).Value = 15
    $ws.Cells.Item($row, 6).Value = 'Fruta'
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = 'Berries'
    $ws.Cells.Item($row, 9).Value = 100112025
    $ws.Cells.Item($row, 10).Value = 'Frutilla'
    $ws.Cells.Item($row, 11).Value = 'Sin especificar'
    $ws.Cells.Item($row, 12).Value = $r[1]
    $ws.Cells.Item($row, 13).Value = $r[2]
    $ws.Cells.Item($row, 14).Value = $r[3]
    $ws.Cells.Item($row, 15).Value = $r[4]
    $ws.Cells.Item($row, 16).Value = $r[5]
    $ws.Cells.Item($row, 17).Value = '$/bandeja 3 kilos'
    $ws.Cells.Item($row, 18).Value = 'Región de Arica y Parinacota'
    $ws.Cells.Item($row, 19).Value = $r[6]
    $ws.Cells.Item($row, 20).Value = 3
}
